$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" message in A1
$ws.Range("A1").Value = "Datos actualizados a 25 de Marzo de 2020 a las 07:16"

# Update country names (column A) per row, reflecting reshuffled country order
$ws.Range("A123").Value = "Mayotte"
$ws.Range("A124").Value = "Honduras"
$ws.Range("A137").Value = "Islas Virgenes de los Estados Unidos"
$ws.Range("A138").Value = "Madagascar"
$ws.Range("A142").Value = "Etiopia"
$ws.Range("A143").Value = "Nueva Caledonia"
$ws.Range("A146").Value = "Guinea Ecuatorial"
$ws.Range("A147").Value = "El Salvador"
$ws.Range("A150").Value = "Haiti"
$ws.Range("A152").Value = "Dominica"
$ws.Range("A157").Value = "Curazao"
$ws.Range("A158").Value = "Islas Caimanes"
$ws.Range("A161").Value = "Guyana"
$ws.Range("A162").Value = "Bahamas"
$ws.Range("A164").Value = "Santa Sede"
$ws.Range("A166").Value = "Guinea"
$ws.Range("A168").Value = "Mozambique"
$ws.Range("A169").Value = "Santa Lucia"
$ws.Range("A170").Value = "Republica del Chad"
$ws.Range("A171").Value = "Birmania"
$ws.Range("A172").Value = "Niger"
$ws.Range("A173").Value = "Angola"
$ws.Range("A174").Value = "Liberia"
$ws.Range("A175").Value = "Zambia"
$ws.Range("A176").Value = "Republica de Africa Central"
$ws.Range("A177").Value = "Republica de Yibuti"
$ws.Range("A178").Value = "Antigua y Barbuda"
$ws.Range("A179").Value = "San Bartolome"
$ws.Range("A181").Value = "Zimbabue"
$ws.Range("A183").Value = "Gambia"
$ws.Range("A184").Value = "Mauritania"
$ws.Range("A185").Value = "Butan"
$ws.Range("A186").Value = "San Martin (Parte Holandesa)"
$ws.Range("A187").Value = "Laos"
$ws.Range("A188").Value = "Nicaragua"
$ws.Range("A190").Value = "Montserrat"
$ws.Range("A191").Value = "Libia"
$ws.Range("A192").Value = "Belice"
$ws.Range("A193").Value = "San Vicente y las Granadinas"
$ws.Range("A194").Value = "Timor Oriental"
$ws.Range("A195").Value = "Granada"
$ws.Range("A196").Value = "Islas Turcas y Caicos"
$ws.Range("A198").Value = "Siria"
$ws.Range("A199").Value = "Somalia"
$ws.Range("A200").Value = "Papua Nueva Guinea"

# Update numeric stats that changed (either real data updates or values that travel with a swapped country)
$ws.Range("B6").Value = 54914
$ws.Range("C6").Value = 33
$ws.Range("D6").Value = 379
$ws.Range("E6").Value = 53751
$ws.Range("H6").Value = 784
$ws.Range("C123").Value = 0
$ws.Range("C124").Value = 6
$ws.Range("C142").Value = 0
$ws.Range("C143").Value = 2
$ws.Range("C146").Value = 0
$ws.Range("C147").Value = 4
$ws.Range("D161").Value = 0
$ws.Range("H161").Value = 1
$ws.Range("D162").Value = 1
$ws.Range("H162").Value = 0
